$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) values

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.974.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5067"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  -2.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07206"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8945"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.874.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07527"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.35%  "

$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008547"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.025.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.031"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.117.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("E23").Value = "  -0.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.422"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.792"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.077"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.56%  "

$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09156"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05147"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7522"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.37%  "

$ws.Range("E36").Value = "  +0.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.223"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.568"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5651"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.86%  "

$ws.Range("E40").Value = "  -1.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.073"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.619"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.531"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4736"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.565"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.04%  "
